$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1793.1428
$ws.Range("I33").Value = 2067.8333
$ws.Range("K33").Value = 2067.8333
$ws.Range("M33").Value = -1838.8333
$ws.Range("H41").Value = 780.4865
$ws.Range("I41").Value = 699.5
$ws.Range("J41").Value = 1032.4445
$ws.Range("K41").Value = 699.5
$ws.Range("L41").Value = 1032.4445
$ws.Range("M41").Value = -259.5
$ws.Range("N41").Value = -1912.4445
$ws.Range("H61").Value = 4252
$ws.Range("I61").Value = 2008
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 6024
$ws.Range("L61").Value = 15000
$ws.Range("M61").Value = -5852
$ws.Range("N61").Value = -15344
$ws.Range("H74").Value = 3939.95
$ws.Range("I74").Value = 3479.9
$ws.Range("J74").Value = 4400
$ws.Range("K74").Value = 3479.9
$ws.Range("L74").Value = 4400
$ws.Range("M74").Value = -2543.9
$ws.Range("N74").Value = -6272
$ws.Range("H77").Value = 3939.95
$ws.Range("I77").Value = 3479.9
$ws.Range("J77").Value = 4400
$ws.Range("K77").Value = 17399.5
$ws.Range("L77").Value = 22000
$ws.Range("M77").Value = -12719.5
$ws.Range("N77").Value = -31360
$ws.Range("H138").Value = 2473.394
$ws.Range("I138").Value = 1801.3077
$ws.Range("K138").Value = 5403.9231
$ws.Range("M138").Value = -263.9231
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 26103.895
$ws.Range("I32").Value = 14946.253
$ws.Range("K32").Value = 14946.253
$ws.Range("M32").Value = -14659.253
$ws.Range("H61").Value = 1290.931
$ws.Range("I61").Value = 880.8261
$ws.Range("J61").Value = 2863
$ws.Range("K61").Value = 880.8261
$ws.Range("L61").Value = 2863
$ws.Range("M61").Value = -668.8261
$ws.Range("N61").Value = -3287
$ws.Range("H92").Value = 52850
$ws.Range("J92").Value = 52850
$ws.Range("L92").Value = 52850
$ws.Range("N92").Value = -57842
$ws.Range("H132").Value = 1834.0682
$ws.Range("I132").Value = 1304.2
$ws.Range("K132").Value = 3912.6
$ws.Range("M132").Value = -1382.6
$ws.Range("H136").Value = 1290.931
$ws.Range("I136").Value = 880.8261
$ws.Range("J136").Value = 2863
$ws.Range("K136").Value = 2642.4783
$ws.Range("L136").Value = 8589
$ws.Range("M136").Value = -92.47829999999976
$ws.Range("N136").Value = -13689
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2142.5
$ws.Range("I99").Value = 1292.1428
$ws.Range("K99").Value = 1292.1428
$ws.Range("M99").Value = 205.8571999999999
$ws.Range("H105").Value = 7905.3076
$ws.Range("I105").Value = 8197
$ws.Range("K105").Value = 8197
$ws.Range("M105").Value = -6450
$ws.Range("H134").Value = 1168.9286
$ws.Range("I134").Value = 1030.5
$ws.Range("K134").Value = 3091.5
$ws.Range("M134").Value = -556.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1661.7368
$ws.Range("I16").Value = 1512.4286
$ws.Range("K16").Value = 1512.4286
$ws.Range("M16").Value = -1225.4286
$ws.Range("H31").Value = 2630.5334
$ws.Range("I31").Value = 1769.091
$ws.Range("J31").Value = 4999.5
$ws.Range("K31").Value = 1769.091
$ws.Range("L31").Value = 4999.5
$ws.Range("M31").Value = -1474.091
$ws.Range("N31").Value = -5589.5
$ws.Range("H34").Value = 2630.5334
$ws.Range("I34").Value = 1769.091
$ws.Range("J34").Value = 4999.5
$ws.Range("K34").Value = 1769.091
$ws.Range("L34").Value = 4999.5
$ws.Range("M34").Value = -1567.091
$ws.Range("N34").Value = -5403.5
$ws.Range("H58").Value = 1980
$ws.Range("I58").Value = 1893.5
$ws.Range("K58").Value = 1893.5
$ws.Range("M58").Value = -1690.5
$ws.Range("H99").Value = 16387.25
$ws.Range("I99").Value = 21925.715
$ws.Range("J99").Value = 8633.4
$ws.Range("K99").Value = 21925.715
$ws.Range("L99").Value = 8633.4
$ws.Range("M99").Value = -20427.715
$ws.Range("N99").Value = -11629.4
$ws.Range("H107").Value = 1367.2963
$ws.Range("I107").Value = 1583.9166
$ws.Range("J107").Value = 1194
$ws.Range("K107").Value = 1583.9166
$ws.Range("L107").Value = 1194
$ws.Range("M107").Value = 336.0834
$ws.Range("N107").Value = -5034
$ws.Range("H113").Value = 1661.7368
$ws.Range("I113").Value = 1512.4286
$ws.Range("K113").Value = 1512.4286
$ws.Range("M113").Value = 657.5714
$ws.Range("H122").Value = 7373.3335
$ws.Range("I122").Value = 7194.4287
$ws.Range("J122").Value = 7999.5
$ws.Range("K122").Value = 21583.2861
$ws.Range("L122").Value = 23998.5
$ws.Range("M122").Value = -19133.2861
$ws.Range("N122").Value = -28898.5
$ws.Range("H126").Value = 16387.25
$ws.Range("I126").Value = 21925.715
$ws.Range("J126").Value = 8633.4
$ws.Range("K126").Value = 65777.145
$ws.Range("L126").Value = 25900.2
$ws.Range("M126").Value = -63307.145
$ws.Range("N126").Value = -30840.2
$ws.Range("H132").Value = 4197.4
$ws.Range("I132").Value = 3996.75
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 11990.25
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -9460.25
$ws.Range("N132").Value = -20060
$ws.Range("H134").Value = 5420.2
$ws.Range("I134").Value = 5420.2
$ws.Range("K134").Value = 16260.6
$ws.Range("M134").Value = -13725.6
$ws.Range("H136").Value = 1980
$ws.Range("I136").Value = 1893.5
$ws.Range("K136").Value = 5680.5
$ws.Range("M136").Value = -3130.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 10629429
$ws.Range("I4").Value = 3753716
$ws.Range("J4").Value = 74802750
$ws.Range("K4").Value = 11261148
$ws.Range("L4").Value = 224408250
$ws.Range("M4").Value = -11261036
$ws.Range("N4").Value = -224408474
$ws.Range("H5").Value = 453.36365
$ws.Range("I5").Value = 453.36365
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1360.09095
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -1248.09095
$ws.Range("N5").ClearContents()
$ws.Range("H124").Value = 9497.5
$ws.Range("I124").Value = 5000
$ws.Range("J124").Value = 10397
$ws.Range("K124").Value = 15000
$ws.Range("L124").Value = 31191
$ws.Range("M124").Value = -10090
$ws.Range("N124").Value = -41011
$ws.Range("H125").Value = 8623.666999999999
$ws.Range("I125").Value = 7997.5
$ws.Range("J125").Value = 9876
$ws.Range("K125").Value = 23992.5
$ws.Range("L125").Value = 29628
$ws.Range("M125").Value = -19072.5
$ws.Range("N125").Value = -39468
$ws.Range("H129").Value = 135923.4
$ws.Range("I129").Value = 200838.5
$ws.Range("J129").Value = 6093.2
$ws.Range("K129").Value = 602515.5
$ws.Range("L129").Value = 18279.6
$ws.Range("M129").Value = -597515.5
$ws.Range("N129").Value = -28279.6
$ws.Range("H130").Value = 14990
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H131").Value = 26199
$ws.Range("I131").Value = 14932.286
$ws.Range("K131").Value = 44796.858
$ws.Range("M131").Value = -39756.858
$ws.Range("H135").Value = 453.36365
$ws.Range("I135").Value = 453.36365
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 4080.27285
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -1545.27285
$ws.Range("N135").ClearContents()
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 5005001.5
$ws.Range("I3").Value = 5005001.5
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 5005001.5
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -5004885.5
$ws.Range("N3").ClearContents()
$ws.Range("H10").Value = 24723.5
$ws.Range("J10").Value = 9631.333000000001
$ws.Range("L10").Value = 9631.333000000001
$ws.Range("N10").Value = -9969.333000000001
$ws.Range("H80").Value = 5995.75
$ws.Range("I80").Value = 5994
$ws.Range("J80").Value = 5996.3335
$ws.Range("K80").Value = 5994
$ws.Range("L80").Value = 5996.3335
$ws.Range("M80").Value = -4996
$ws.Range("N80").Value = -7992.3335
$ws.Range("H83").Value = 5995.75
$ws.Range("I83").Value = 5994
$ws.Range("J83").Value = 5996.3335
$ws.Range("K83").Value = 29970
$ws.Range("L83").Value = 29981.6675
$ws.Range("M83").Value = -24978
$ws.Range("N83").Value = -39965.6675
$ws.Range("H102").Value = 2733.6
$ws.Range("I102").Value = 2802.75
$ws.Range("K102").Value = 2802.75
$ws.Range("M102").Value = -1180.75
$ws.Range("H122").Value = 3546.5334
$ws.Range("I122").Value = 3546.5334
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 10639.6002
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -8189.600199999999
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 3415.647
$ws.Range("I126").Value = 3555.75
$ws.Range("K126").Value = 10667.25
$ws.Range("M126").Value = -8197.25
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1645.6
$ws.Range("J22").Value = 1699.5
$ws.Range("L22").Value = 1699.5
$ws.Range("N22").Value = -2289.5
$ws.Range("H27").Value = 1645.6
$ws.Range("J27").Value = 1699.5
$ws.Range("L27").Value = 1699.5
$ws.Range("N27").Value = -1913.5
$ws.Range("H46").Value = 26818.412
$ws.Range("J46").Value = 1949.625
$ws.Range("L46").Value = 1949.625
$ws.Range("N46").Value = -2325.625
$ws.Range("H100").Value = 23446.65
$ws.Range("I100").Value = 8371.666999999999
$ws.Range("K100").Value = 8371.666999999999
$ws.Range("M100").Value = -7830.666999999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 4999.5
$ws.Range("J7").Value = 4999.5
$ws.Range("L7").Value = 4999.5
$ws.Range("N7").Value = -5225.5
$ws.Range("H136").Value = 297
$ws.Range("I136").Value = 297
$ws.Range("K136").Value = 891
$ws.Range("M136").Value = 1659
